$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim() -eq "") {
        $p.Range.Delete()
        break
    }
}
